$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose new values look numeric need to be forced to Text so Excel
# keeps the literal string (matching trailing zeros / precision in the source)
# instead of auto-converting them into floating point numbers.
$numericTextCells = @("D5","D7","D8","D9","D10","D11","D13","D15","D17","D18","D19","D20","D22","D26","D27","D29","D30","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.456.20"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.827.70"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "315.37"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "0.5141"
$ws.Range("E7").Value = "  -3.61%  "
$ws.Range("D8").Value = "0.3926"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").Value = "0.07668"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.110"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "41.67"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "6.288"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "7.541"
$ws.Range("D16").Value = "1.824.06"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "93.57"
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("D18").Value = "0.00001106"
$ws.Range("E18").Value = "  +3.06%  "
$ws.Range("D19").Value = "0.06674"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").Value = "17.66"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "6.146"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("D23").Value = "28.487.47"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  +6.97%  "
$ws.Range("D26").Value = "20.81"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").Value = "156.94"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "2.037.55"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").Value = "2.391"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").Value = "124.49"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "0.1085"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("D33").Value = "5.652"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").Value = "3.664"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").Value = "0.07031"
$ws.Range("E35").Value = "  -4.59%  "
$ws.Range("D36").Value = "0.2209"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("D37").Value = "8.914"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").Value = "0.02324"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "5.155"
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("D40").Value = "0.6261"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "11.20"
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("D42").Value = "1.177"
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "1.393"
$ws.Range("E44").Value = "  -1.58%  "
$ws.Range("D45").Value = "13.43"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "0.5903"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "3.710"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "125.21"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").Value = "1.975"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").Value = "1.197"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "0.06929"
$ws.Range("E51").Value = "  +0.60%  "
